$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024878486356477
$ws.Cells.Item(2, 4).Value = 1.030165460717941
$ws.Cells.Item(2, 5).Value = 1.035649366270501
$ws.Cells.Item(2, 6).Value = 1.04753693780185
$ws.Cells.Item(2, 9).Value = 1.032560302146778
$ws.Cells.Item(2, 10).Value = 1.030050636249559
$ws.Cells.Item(2, 11).Value = 1.032977262253909
$ws.Cells.Item(2, 12).Value = 1.038445365765775
$ws.Cells.Item(2, 13).Value = 1.050299299224731
$ws.Cells.Item(2, 14).Value = 1.01406302036066
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025752978934303
$ws.Cells.Item(3, 4).Value = 1.030816392315695
$ws.Cells.Item(3, 5).Value = 1.03645691777828
$ws.Cells.Item(3, 6).Value = 1.048502071177893
$ws.Cells.Item(3, 9).Value = 1.032723270124941
$ws.Cells.Item(3, 10).Value = 1.030564611160075
$ws.Cells.Item(3, 11).Value = 1.033436914725521
$ws.Cells.Item(3, 12).Value = 1.039062360842412
$ws.Cells.Item(3, 13).Value = 1.051075899308887
$ws.Cells.Item(3, 14).Value = 1.014233418889214
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.026319157645383
$ws.Cells.Item(4, 4).Value = 1.031237414278182
$ws.Cells.Item(4, 5).Value = 1.036980129289829
$ws.Cells.Item(4, 6).Value = 1.0491273180001
$ws.Cells.Item(4, 9).Value = 1.032826980742789
$ws.Cells.Item(4, 10).Value = 1.030896884142857
$ws.Cells.Item(4, 11).Value = 1.033733483033454
$ws.Cells.Item(4, 12).Value = 1.039461617865554
$ws.Cells.Item(4, 13).Value = 1.051578532144539
$ws.Cells.Item(4, 14).Value = 1.014343553498697
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.026557255352977
$ws.Cells.Item(5, 4).Value = 1.031414368527126
$ws.Cells.Item(5, 5).Value = 1.037200246681133
$ws.Cells.Item(5, 6).Value = 1.04939034760694
$ws.Cells.Item(5, 9).Value = 1.032870163177535
$ws.Cells.Item(5, 10).Value = 1.031036498024277
$ws.Cells.Item(5, 11).Value = 1.033857953791685
$ws.Cells.Item(5, 12).Value = 1.039629469066817
$ws.Cells.Item(5, 13).Value = 1.051789866460172
$ws.Cells.Item(5, 14).Value = 1.014389823853464
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026597237435687
$ws.Cells.Item(6, 4).Value = 1.031444077331858
$ws.Cells.Item(6, 5).Value = 1.037237214654893
$ws.Cells.Item(6, 6).Value = 1.049434521701664
$ws.Cells.Item(6, 9).Value = 1.03287738919164
$ws.Cells.Item(6, 10).Value = 1.031059935456446
$ws.Cells.Item(6, 11).Value = 1.033878840830481
$ws.Cells.Item(6, 12).Value = 1.039657652194444
$ws.Cells.Item(6, 13).Value = 1.051825351991506
$ws.Cells.Item(6, 14).Value = 1.014397591065044
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026322338820883
$ws.Cells.Item(7, 4).Value = 1.031239778922079
$ws.Cells.Item(7, 5).Value = 1.036983069885286
$ws.Cells.Item(7, 6).Value = 1.049130831925593
$ws.Cells.Item(7, 9).Value = 1.03282755939015
$ws.Cells.Item(7, 10).Value = 1.030898749961282
$ws.Cells.Item(7, 11).Value = 1.033735147031035
$ws.Cells.Item(7, 12).Value = 1.039463860689124
$ws.Cells.Item(7, 13).Value = 1.051581355897845
$ws.Cells.Item(7, 14).Value = 1.014344171884747
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025173957458756
$ws.Cells.Item(8, 4).Value = 1.030385481198091
$ws.Cells.Item(8, 5).Value = 1.035922142086748
$ws.Cells.Item(8, 6).Value = 1.047862954785455
$ws.Cells.Item(8, 9).Value = 1.032615737843937
$ws.Cells.Item(8, 10).Value = 1.030224398583366
$ws.Cells.Item(8, 11).Value = 1.033132780774336
$ws.Cells.Item(8, 12).Value = 1.03865387727248
$ws.Cells.Item(8, 13).Value = 1.050561729029075
$ws.Cells.Item(8, 14).Value = 1.014120632908472
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.023152896975951
$ws.Cells.Item(9, 4).Value = 1.028878837222609
$ws.Cells.Item(9, 5).Value = 1.034057863530288
$ws.Cells.Item(9, 6).Value = 1.045634536745269
$ws.Cells.Item(9, 9).Value = 1.032229186365226
$ws.Cells.Item(9, 10).Value = 1.029033832124334
$ws.Cells.Item(9, 11).Value = 1.032064821493781
$ws.Cells.Item(9, 12).Value = 1.03722678985006
$ws.Cells.Item(9, 13).Value = 1.048766002451267
$ws.Cells.Item(9, 14).Value = 1.013725791633527
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021807301920458
$ws.Cells.Item(10, 4).Value = 1.027873656628635
$ws.Cells.Item(10, 5).Value = 1.032818605497332
$ws.Cells.Item(10, 6).Value = 1.044152877419586
$ws.Cells.Item(10, 9).Value = 1.031962594331995
$ws.Cells.Item(10, 10).Value = 1.028238657662767
$ws.Cells.Item(10, 11).Value = 1.031348543162167
$ws.Cells.Item(10, 12).Value = 1.036275609461685
$ws.Cells.Item(10, 13).Value = 1.047569597407261
$ws.Cells.Item(10, 14).Value = 1.013461957485844
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.021225081662648
$ws.Cells.Item(11, 4).Value = 1.027438244463909
$ws.Cells.Item(11, 5).Value = 1.032282864430456
$ws.Cells.Item(11, 6).Value = 1.043512258988037
$ws.Cells.Item(11, 9).Value = 1.031845058058781
$ws.Cells.Item(11, 10).Value = 1.027894004684867
$ws.Cells.Item(11, 11).Value = 1.031037380798006
$ws.Cells.Item(11, 12).Value = 1.035863802947705
$ws.Cells.Item(11, 13).Value = 1.047051734950101
$ws.Cells.Item(11, 14).Value = 1.013347575365701
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.021008885128114
$ws.Cells.Item(12, 4).Value = 1.027276490224667
$ws.Cells.Item(12, 5).Value = 1.032083997924398
$ws.Cells.Item(12, 6).Value = 1.043274448969648
$ws.Cells.Item(12, 9).Value = 1.031801085083161
$ws.Cells.Item(12, 10).Value = 1.027765935621124
$ws.Cells.Item(12, 11).Value = 1.030921650828703
$ws.Cells.Item(12, 12).Value = 1.035710850000278
$ws.Cells.Item(12, 13).Value = 1.046859407495104
$ws.Cells.Item(12, 14).Value = 1.013305068073951
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.021055257013138
$ws.Cells.Item(13, 4).Value = 1.027311188062821
$ws.Cells.Item(13, 5).Value = 1.032126649486733
$ws.Cells.Item(13, 6).Value = 1.043325453477042
$ws.Cells.Item(13, 9).Value = 1.031810531667634
$ws.Cells.Item(13, 10).Value = 1.0277934090977
$ws.Cells.Item(13, 11).Value = 1.030946482087315
$ws.Cells.Item(13, 12).Value = 1.035743658440412
$ws.Cells.Item(13, 13).Value = 1.046900661020538
$ws.Cells.Item(13, 14).Value = 1.013314186961978
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.021207209436505
$ws.Cells.Item(14, 4).Value = 1.027424874263317
$ws.Cells.Item(14, 5).Value = 1.032266423370187
$ws.Cells.Item(14, 6).Value = 1.043492598587683
$ws.Cells.Item(14, 9).Value = 1.031841429655763
$ws.Cells.Item(14, 10).Value = 1.027883419461873
$ws.Cells.Item(14, 11).Value = 1.031027817584116
$ws.Cells.Item(14, 12).Value = 1.03585115959185
$ws.Cells.Item(14, 13).Value = 1.047035836488602
$ws.Cells.Item(14, 14).Value = 1.013344062119349
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.021300841120086
$ws.Cells.Item(15, 4).Value = 1.027494917113882
$ws.Cells.Item(15, 5).Value = 1.032352560153274
$ws.Cells.Item(15, 6).Value = 1.043595601357995
$ws.Cells.Item(15, 9).Value = 1.031860425238412
$ws.Cells.Item(15, 10).Value = 1.027938871276759
$ws.Cells.Item(15, 11).Value = 1.031077911172649
$ws.Cells.Item(15, 12).Value = 1.035917396003059
$ws.Cells.Item(15, 13).Value = 1.047119126534316
$ws.Cells.Item(15, 14).Value = 1.013362466461016
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.021845951150279
$ws.Cells.Item(16, 4).Value = 1.027902550217774
$ws.Cells.Item(16, 5).Value = 1.032854179237657
$ws.Cells.Item(16, 6).Value = 1.044195413294094
$ws.Cells.Item(16, 9).Value = 1.031970350643875
$ws.Cells.Item(16, 10).Value = 1.02826152413191
$ws.Cells.Item(16, 11).Value = 1.031369172824703
$ws.Cells.Item(16, 12).Value = 1.036302941087234
$ws.Cells.Item(16, 13).Value = 1.047603970326643
$ws.Cells.Item(16, 14).Value = 1.013469545729457
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02218800058213
$ws.Cells.Item(17, 4).Value = 1.028158205374682
$ws.Cells.Item(17, 5).Value = 1.033169064504662
$ws.Cells.Item(17, 6).Value = 1.044571914823926
$ws.Cells.Item(17, 9).Value = 1.032038742207465
$ws.Cells.Item(17, 10).Value = 1.02846382621551
$ws.Cells.Item(17, 11).Value = 1.031551604239957
$ws.Cells.Item(17, 12).Value = 1.036544800465034
$ws.Cells.Item(17, 13).Value = 1.04790815153572
$ws.Cells.Item(17, 14).Value = 1.0135366764626
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022387553760764
$ws.Cells.Item(18, 4).Value = 1.028307308931587
$ws.Cells.Item(18, 5).Value = 1.033352815194172
$ws.Cells.Item(18, 6).Value = 1.044791613341843
$ws.Cells.Item(18, 9).Value = 1.032078431163617
$ws.Cells.Item(18, 10).Value = 1.02858179299368
$ws.Cells.Item(18, 11).Value = 1.031657915952619
$ws.Cells.Item(18, 12).Value = 1.036685878742203
$ws.Cells.Item(18, 13).Value = 1.048085593416751
$ws.Cells.Item(18, 14).Value = 1.013575819131258
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022455603245883
$ws.Cells.Item(19, 4).Value = 1.028358146667092
$ws.Cells.Item(19, 5).Value = 1.033415483538131
$ws.Cells.Item(19, 6).Value = 1.04486654034131
$ws.Cells.Item(19, 9).Value = 1.03209192966995
$ws.Cells.Item(19, 10).Value = 1.028622011054016
$ws.Cells.Item(19, 11).Value = 1.031694148934407
$ws.Cells.Item(19, 12).Value = 1.036733983735684
$ws.Cells.Item(19, 13).Value = 1.048146099546377
$ws.Cells.Item(19, 14).Value = 1.013589163459373
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.022151297610676
$ws.Cells.Item(20, 4).Value = 1.0281307776154
$ws.Cells.Item(20, 5).Value = 1.033135271654054
$ws.Cells.Item(20, 6).Value = 1.044531510293184
$ws.Cells.Item(20, 9).Value = 1.032031425399207
$ws.Cells.Item(20, 10).Value = 1.028442124481289
$ws.Cells.Item(20, 11).Value = 1.031532041146176
$ws.Cells.Item(20, 12).Value = 1.036518850634861
$ws.Cells.Item(20, 13).Value = 1.047875513885861
$ws.Cells.Item(20, 14).Value = 1.013529475369249
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.021162461403802
$ws.Cells.Item(21, 4).Value = 1.027391397126582
$ws.Cells.Item(21, 5).Value = 1.032225259801018
$ws.Cells.Item(21, 6).Value = 1.043443374526265
$ws.Cells.Item(21, 9).Value = 1.031832339651526
$ws.Cells.Item(21, 10).Value = 1.02785691501386
$ws.Cells.Item(21, 11).Value = 1.031003870450088
$ws.Cells.Item(21, 12).Value = 1.035819502887573
$ws.Cells.Item(21, 13).Value = 1.046996029851176
$ws.Cells.Item(21, 14).Value = 1.013335265197907
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.020541123239411
$ws.Cells.Item(22, 4).Value = 1.026926388295065
$ws.Cells.Item(22, 5).Value = 1.031653861251306
$ws.Cells.Item(22, 6).Value = 1.042760055359949
$ws.Cells.Item(22, 9).Value = 1.031705345747107
$ws.Cells.Item(22, 10).Value = 1.027488684133546
$ws.Cells.Item(22, 11).Value = 1.030670918917207
$ws.Cells.Item(22, 12).Value = 1.035379855537276
$ws.Cells.Item(22, 13).Value = 1.046443235946552
$ws.Cells.Item(22, 14).Value = 1.013213038098742
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.020870469708917
$ws.Cells.Item(23, 4).Value = 1.027172910152326
$ws.Cells.Item(23, 5).Value = 1.031956697703807
$ws.Cells.Item(23, 6).Value = 1.04312221613585
$ws.Cells.Item(23, 9).Value = 1.03177283995744
$ws.Cells.Item(23, 10).Value = 1.027683917091474
$ws.Cells.Item(23, 11).Value = 1.03084750485942
$ws.Cells.Item(23, 12).Value = 1.03561291480285
$ws.Cells.Item(23, 13).Value = 1.046736265642195
$ws.Cells.Item(23, 14).Value = 1.013277844192121
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.022167881975386
$ws.Cells.Item(24, 4).Value = 1.028143171088108
$ws.Cells.Item(24, 5).Value = 1.033150540930987
$ws.Cells.Item(24, 6).Value = 1.04454976707871
$ws.Cells.Item(24, 9).Value = 1.032034732176242
$ws.Cells.Item(24, 10).Value = 1.028451930661687
$ws.Cells.Item(24, 11).Value = 1.031540881167379
$ws.Cells.Item(24, 12).Value = 1.036530576228331
$ws.Cells.Item(24, 13).Value = 1.047890261379072
$ws.Cells.Item(24, 14).Value = 1.013532729275279
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02367508153529
$ws.Cells.Item(25, 4).Value = 1.029268480324504
$ws.Cells.Item(25, 5).Value = 1.03453919729021
$ws.Cells.Item(25, 6).Value = 1.046209946195543
$ws.Cells.Item(25, 9).Value = 1.032330689771094
$ws.Cells.Item(25, 10).Value = 1.029341884778403
$ws.Cells.Item(25, 11).Value = 1.032341678796014
$ws.Cells.Item(25, 12).Value = 1.037595694197251
$ws.Cells.Item(25, 13).Value = 1.049230115386164
$ws.Cells.Item(25, 14).Value = 1.013827976197323
